$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Linhas")

# --- Recolor existing bus-line rows ---
# 600-604 (Autocarro): #1d4ed8 -> #f94144
$ws.Range("B3:B7").Value = "#f94144"
# 700-704 (Autocarro): #1d4ed8 -> #f3722c
$ws.Range("B8:B12").Value = "#f3722c"

# --- Make room for the new rows (705, 706, 707) above the old trailing row 14 ---
# Row 14 (old) was an otherwise-empty row with a residual style on B14; pushing
# it down to row 16 reproduces that leftover formatting ending up on the new B16 cell.
$ws.Range("A14:A15").EntireRow.Insert()

# --- New row 13: 705 (Autocarro) ---
$ws.Range("D13").Value = "HSJ2,RCA,FEL1,CVD1,ISCP1,DRB1,OG1,MOD1,SMM4,LARJ1,THB1,LNC1,PRD3,RGLE2,GGC2,FMGO3,AS1,RMOS1,AST1,CRG1,LVDR1,JLD2,MEIL1,MRT1,GAN1,CML1,MSHP3,RIGR1,ERMI1,ERM2,ERM4,VASG2,RTEL2,SA4,QCH2,FMGC1,CHAO1,MCST2,9AG1,COVN2,FND2,4CAM3,TAZ2,FS1,ASRR3,RCMB1,FTSR1,CFER1,25AB1,PAV1,CLVR1,CMV,VALC1,ILHA,PRS3,RBQ1,FTMR1,EMDN1,PRA1,JLS1,VLG5"
$ws.Range("A13").Value = "705 (Autocarro)"
$ws.Range("B13").Value = "#f3722c"
$ws.Range("C13").Value = 3

# --- New row 14: 706 (Autocarro) ---
$ws.Range("D14").Value = "HSJ1,ASP3,ENX1,ARSM3,ARSF3,HUMB1,TNG1,PIV1,GLV1,SGMC1,DDN1,5OUT1,DCPF1,SRG1,FMGO1,GGC1,RGLE1,PRD1,TRM1,RGLO1,PMAI1,AGR1,VESS1,RARR1,RCNH1,CSAD1,ANV1,ZIND1,MPN1,PRMV1,PNHL1,TARG1,ARDG3,PHIP1,SAMP1,TSL1,TVG1,RLC2,RGFT2,RIGR1,ERMI1,ERM2"
$ws.Range("A14").Value = "706 (Autocarro)"
$ws.Range("B14").Value = "#f3722c"
$ws.Range("C14").Value = 3

# --- New row 15: 707 (Autocarro) ---
$ws.Range("D15").Value = "ERM2,ERM4,ERMI2,RIGR2,RGFT1,RLC1,TVG2,TSL2,SAMP2,PHIP2,ARDG2,RCTL2,PISA2,ARRV2,QULH2,ARRG2,IGMI2,TESC2,CSAD4,RCNH2,RARR2,VESS2,AGR2,PMAI2,RGLO2,TRM2,PRD2,RGLE2,GGC2,FMGO2,SRG2,DCPF2,5OUT2,DDN2,SGMC2,GLV2,PIV2,TNG2,HUMB2,ARSF4,ARSM2,ARS5,ENX2,ASP4,HSJ10,HSJ1"
$ws.Range("A15").Value = "707 (Autocarro)"
$ws.Range("B15").Value = "#f3722c"
$ws.Range("C15").Value = 3

# --- Row 16 (pre-existing, previously-empty B16 now filled in): 800 (Autocarro) ---
$ws.Range("D16").Value = "SR2,TPT1,SCT1,STBX1,OTH1,INFS1,ASOT1,DAC1,CV1,CB1,VARZ1,STO1,STCZ1,STE1,CMR1,VL1,TAR5,VINH1,ICPU1,25A1,CGG,GODC1,BVG1,MTCR1,1RI1,1GFE1,GODS"
$ws.Range("A16").Value = "800 (Autocarro)"
$ws.Range("C16").Value = 3
$ws.Range("B16").Value = "#577590"

# --- Selection matches the author's last-saved cursor position ---
$ws.Range("B16").Select()
